$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '55.404.07'
$ws.Range('E2').Value = '  -2.87%  '
$ws.Range('D3').Value = '2.942.05'
$ws.Range('E3').Value = '  -5.09%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '488.15'
$ws.Range('E5').Value = '  -5.82%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '131.55'
$ws.Range('E6').Value = '  -1.38%  '
$ws.Range('E7').Value = '  -0.13%  '
$ws.Range('D8').Value = '2.939.53'
$ws.Range('E8').Value = '  -5.00%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.420'
$ws.Range('E9').Value = '  -5.25%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '7.08'
$ws.Range('E10').Value = '  -0.30%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.101'
$ws.Range('E11').Value = '  -6.34%  '
$ws.Range('E12').Value = '  -8.59%  '
$ws.Range('E13').Value = '  +0.63%  '
$ws.Range('D14').Value = '3.443.38'
$ws.Range('E14').Value = '  -5.69%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '24.50'
$ws.Range('E15').Value = '  -2.83%  '
$ws.Range('D16').Value = '55.427.71'
$ws.Range('E16').Value = '  -2.90%  '
$ws.Range('D17').Value = '2.950.11'
$ws.Range('E17').Value = '  -5.11%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.0000140'
$ws.Range('E18').Value = '  -5.56%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '5.58'
$ws.Range('E19').Value = '  -2.23%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.09'
$ws.Range('E20').Value = '  -5.92%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '7.49'
$ws.Range('E21').Value = '  -5.28%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '316.05'
$ws.Range('E22').Value = '  -7.40%  '
$ws.Range('E23').Value = '  -0.03%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.458'
$ws.Range('E24').Value = '  -8.68%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '59.90'
$ws.Range('E25').Value = '  -12.08%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.01'
$ws.Range('E26').Value = '  +0.70%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.159'
$ws.Range('E27').Value = '  -3.51%  '
$ws.Range('E28').Value = '  +0.06%  '
$ws.Range('D29').Value = '0.0₃0833'
$ws.Range('E29').Value = '  -8.96%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.42'
$ws.Range('E30').Value = '  -3.58%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.52'
$ws.Range('E31').Value = '  -5.86%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.15'
$ws.Range('E32').Value = '  -4.40%  '
$ws.Range('E33').Value = '  -8.54%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '19.36'
$ws.Range('E34').Value = '  -9.29%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '147.46'
$ws.Range('E35').Value = '  -5.83%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.35'
$ws.Range('E36').Value = '  -8.21%  '
$ws.Range('E37').Value = '  -5.49%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '5.65'
$ws.Range('E38').Value = '  -7.48%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0651'
$ws.Range('E39').Value = '  -4.46%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '22.96'
$ws.Range('E40').Value = '  -8.57%  '
$ws.Range('D41').Value = '2.970.64'
$ws.Range('E41').Value = '  -5.29%  '
$ws.Range('E42').Value = '  +0.12%  '
$ws.Range('E43').Value = '  -9.99%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.983'
$ws.Range('E44').Value = '  -5.64%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.625'
$ws.Range('E45').Value = '  -7.47%  '
$ws.Range('E46').Value = '  -4.87%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.51'
$ws.Range('E47').Value = '  -9.38%  '
$ws.Range('D48').Value = '2.106.16'
$ws.Range('E48').Value = '  -6.47%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0233'
$ws.Range('E49').Value = '  +0.66%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '19.18'
$ws.Range('E50').Value = '  -2.56%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '5.52'
$ws.Range('E51').Value = '  -9.83%  '
